$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an untouched data cell (default style, no number formatting)
# to restore style after assignment, since Excel auto-applies a
# "quote prefix" / text style to cells when values are entered that
# look like numbers but should remain literal text.
$refStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'55.608.25"
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = "'  +5.66%  "
$ws.Range("E2").Style = $refStyle

$ws.Range("D3").Value = "'2.517.62"
$ws.Range("D3").Style = $refStyle
$ws.Range("E3").Value = "'  +7.63%  "
$ws.Range("E3").Style = $refStyle

$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = $refStyle

$ws.Range("D5").Value = "'492.39"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "'  +11.94%  "
$ws.Range("E5").Style = $refStyle

$ws.Range("D6").Value = "'141.39"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "'  +14.77%  "
$ws.Range("E6").Style = $refStyle

$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E7").Style = $refStyle

$ws.Range("D8").Value = "'0.516"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "'  +9.21%  "
$ws.Range("E8").Style = $refStyle

$ws.Range("D9").Value = "'2.513.65"
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = "'  +7.17%  "
$ws.Range("E9").Style = $refStyle

$ws.Range("D10").Value = "'0.0994"
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = "'  +10.73%  "
$ws.Range("E10").Style = $refStyle

$ws.Range("D11").Value = "'5.55"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "'  +5.70%  "
$ws.Range("E11").Style = $refStyle

$ws.Range("D12").Value = "'0.333"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "'  +8.09%  "
$ws.Range("E12").Style = $refStyle

$ws.Range("E13").Value = "'  +1.96%  "
$ws.Range("E13").Style = $refStyle

$ws.Range("D14").Value = "'2.917.07"
$ws.Range("D14").Style = $refStyle
$ws.Range("E14").Value = "'  +6.91%  "
$ws.Range("E14").Style = $refStyle

$ws.Range("D15").Value = "'55.630.23"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "'  +5.73%  "
$ws.Range("E15").Style = $refStyle

$ws.Range("D16").Value = "'20.96"
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = "'  +9.39%  "
$ws.Range("E16").Style = $refStyle

$ws.Range("D17").Value = "'0.0000139"
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = "'  +15.61%  "
$ws.Range("E17").Style = $refStyle

$ws.Range("D18").Value = "'2.506.37"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "'  +6.97%  "
$ws.Range("E18").Style = $refStyle

$ws.Range("D19").Value = "'4.42"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "'  +9.56%  "
$ws.Range("E19").Style = $refStyle

$ws.Range("D20").Value = "'323.89"
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = "'  +6.48%  "
$ws.Range("E20").Style = $refStyle

$ws.Range("D21").Value = "'10.08"
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = "'  +10.81%  "
$ws.Range("E21").Style = $refStyle

$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "'  -0.15%  "
$ws.Range("E22").Style = $refStyle

$ws.Range("D23").Value = "'5.79"
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = "'  +10.43%  "
$ws.Range("E23").Style = $refStyle

$ws.Range("D24").Value = "'58.36"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "'  +7.81%  "
$ws.Range("E24").Style = $refStyle

$ws.Range("D25").Value = "'0.172"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "'  +13.48%  "
$ws.Range("E25").Style = $refStyle

$ws.Range("D26").Value = "'0.414"
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = "'  +11.48%  "
$ws.Range("E26").Style = $refStyle

$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "'  +0.41%  "
$ws.Range("E27").Style = $refStyle

$ws.Range("D28").Value = "'2.614.75"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "'  +7.93%  "
$ws.Range("E28").Style = $refStyle

$ws.Range("D29").Value = "'7.49"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "'  +4.51%  "
$ws.Range("E29").Style = $refStyle

$ws.Range("D30").Value = "'0.0₃0803"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "'  +16.82%  "
$ws.Range("E30").Style = $refStyle

$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = "'  +0.31%  "
$ws.Range("E31").Style = $refStyle

$ws.Range("D32").Value = "'150.66"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "'  +4.34%  "
$ws.Range("E32").Style = $refStyle

$ws.Range("D33").Value = "'18.35"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "'  +6.72%  "
$ws.Range("E33").Style = $refStyle

$ws.Range("D34").Value = "'1.51"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "'  +12.16%  "
$ws.Range("E34").Style = $refStyle

$ws.Range("D35").Value = "'5.25"
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = "'  +8.83%  "
$ws.Range("E35").Style = $refStyle

$ws.Range("D36").Value = "'0.880"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "'  +4.85%  "
$ws.Range("E36").Style = $refStyle

$ws.Range("B37").Value = "'ImmutableX"
$ws.Range("B37").Style = $refStyle
$ws.Range("C37").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C37").Style = $refStyle
$ws.Range("D37").Value = "'1.14"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "'  +12.29%  "
$ws.Range("E37").Style = $refStyle

$ws.Range("B38").Value = "'NEARProtocol"
$ws.Range("B38").Style = $refStyle
$ws.Range("C38").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C38").Style = $refStyle
$ws.Range("D38").Value = "'3.73"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "'  +4.99%  "
$ws.Range("E38").Style = $refStyle

$ws.Range("D39").Value = "'34.52"
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = "'  +7.10%  "
$ws.Range("E39").Style = $refStyle

$ws.Range("D40").Value = "'0.616"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "'  +16.48%  "
$ws.Range("E40").Style = $refStyle

$ws.Range("E41").Value = "'  +9.58%  "
$ws.Range("E41").Style = $refStyle

$ws.Range("D42").Value = "'0.995"
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = "'  +0.12%  "
$ws.Range("E42").Style = $refStyle

$ws.Range("D43").Value = "'3.46"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "'  +7.91%  "
$ws.Range("E43").Style = $refStyle

$ws.Range("E44").Value = "'  +7.86%  "
$ws.Range("E44").Style = $refStyle

$ws.Range("D45").Value = "'4.79"
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = "'  +17.25%  "
$ws.Range("E45").Style = $refStyle

$ws.Range("D46").Value = "'2.008.62"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "'  +3.97%  "
$ws.Range("E46").Style = $refStyle

$ws.Range("D47").Value = "'259.23"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "'  +34.02%  "
$ws.Range("E47").Style = $refStyle

$ws.Range("D48").Value = "'0.0917"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "'  +9.64%  "
$ws.Range("E48").Style = $refStyle

$ws.Range("E49").Value = "'  -0.15%  "
$ws.Range("E49").Style = $refStyle

$ws.Range("D50").Value = "'0.0227"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "'  +8.02%  "
$ws.Range("E50").Style = $refStyle

$ws.Range("D51").Value = "'17.60"
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = "'  +10.88%  "
$ws.Range("E51").Style = $refStyle
